$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 109.026058
$ws.Range("H2").Value = 327.078174
$ws.Range("I2").Value = 0.3049840938689738
$ws.Range("J2").Value = 0.3049840938689738
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 3714.599165841286
$ws.Range("R2").Value = 33431.39249257158
$ws.Range("S2").Value = 0.1577214413026566
$ws.Range("T2").Value = 0.1577214413026567

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 109.026058
$ws.Range("H3").Value = 327.078174
$ws.Range("I3").Value = 0.3049840938689738
$ws.Range("J3").Value = 0.3049840938689738
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.056491
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 3018.440601680826
$ws.Range("R3").Value = 27165.96541512743
$ws.Range("S3").Value = 0.1281626309943287
$ws.Range("T3").Value = 0.1281626309943287

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 109.026058
$ws.Range("H4").Value = 327.078174
$ws.Range("I4").Value = 0.3049840938689738
$ws.Range("J4").Value = 0.3049840938689738
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 449.8368998715446
$ws.Range("R4").Value = 4048.532098843902
$ws.Range("S4").Value = 0.01910002157198845
$ws.Range("T4").Value = 0.01910002157198845

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 89.97721833333333
$ws.Range("H5").Value = 269.931655
$ws.Range("I5").Value = 0.2516978134001918
$ws.Range("J5").Value = 0.2516978134001917
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 3065.590981613948
$ws.Range("R5").Value = 27590.31883452554
$ws.Range("S5").Value = 0.1301646305504061
$ws.Range("T5").Value = 0.1301646305504061

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 89.97721833333333
$ws.Range("H6").Value = 269.931655
$ws.Range("I6").Value = 0.2516978134001918
$ws.Range("J6").Value = 0.2516978134001917
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.056491
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 2491.064008235845
$ws.Range("R6").Value = 22419.5760741226
$ws.Range("S6").Value = 0.1057702832028573
$ws.Range("T6").Value = 0.1057702832028573

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 89.97721833333333
$ws.Range("H7").Value = 269.931655
$ws.Range("I7").Value = 0.2516978134001918
$ws.Range("J7").Value = 0.2516978134001917
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 371.2421938077572
$ws.Range("R7").Value = 3341.179744269814
$ws.Range("S7").Value = 0.01576289964692827
$ws.Range("T7").Value = 0.01576289964692827

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 158.477852
$ws.Range("H8").Value = 475.433556
$ws.Range("I8").Value = 0.4433180927308344
$ws.Range("J8").Value = 0.4433180927308344
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 5399.458694943541
$ws.Range("R8").Value = 48595.12825449187
$ws.Range("S8").Value = 0.2292603776611744
$ws.Range("T8").Value = 0.2292603776611744

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 158.477852
$ws.Range("H9").Value = 475.433556
$ws.Range("I9").Value = 0.4433180927308344
$ws.Range("J9").Value = 0.4433180927308344
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.056491
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 4387.538096112443
$ws.Range("R9").Value = 39487.84286501199
$ws.Range("S9").Value = 0.1862943487019391
$ws.Range("T9").Value = 0.1862943487019391

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 158.477852
$ws.Range("H10").Value = 475.433556
$ws.Range("I10").Value = 0.4433180927308344
$ws.Range("J10").Value = 0.4433180927308344
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 653.8729084562652
$ws.Range("R10").Value = 5884.856176106387
$ws.Range("S10").Value = 0.02776336636772094
$ws.Range("T10").Value = 0.02776336636772093

